$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1236.5625
$ws.Range("I17").Value = 660
$ws.Range("J17").Value = 1428.75
$ws.Range("K17").Value = 1980
$ws.Range("L17").Value = 4286.25
$ws.Range("M17").Value = -1812
$ws.Range("N17").Value = -4622.25

$ws.Range("H19").Value = 8724.25
$ws.Range("I19").Value = 7499
$ws.Range("J19").Value = 9949.5
$ws.Range("K19").Value = 7499
$ws.Range("L19").Value = 9949.5
$ws.Range("M19").Value = -7324
$ws.Range("N19").Value = -10299.5

$ws.Range("H33").Value = 1104.2142
$ws.Range("I33").Value = 1150.6923
$ws.Range("K33").Value = 1150.6923
$ws.Range("M33").Value = -921.6922999999999

$ws.Range("H40").Value = 3418
$ws.Range("I40").Value = 2933.6365
$ws.Range("J40").Value = 3773.2
$ws.Range("K40").Value = 2933.6365
$ws.Range("L40").Value = 3773.2
$ws.Range("M40").Value = -2758.6365
$ws.Range("N40").Value = -4123.2

$ws.Range("H64").Value = 31256882
$ws.Range("J64").Value = 7333.3335
$ws.Range("L64").Value = 7333.3335
$ws.Range("N64").Value = -7829.3335

$ws.Range("H67").Value = 31256882
$ws.Range("J67").Value = 7333.3335
$ws.Range("L67").Value = 7333.3335
$ws.Range("N67").Value = -9049.333500000001

$ws.Range("H70").Value = 5427.375
$ws.Range("J70").Value = 3651.842
$ws.Range("L70").Value = 10955.526
$ws.Range("N70").Value = -11495.526

$ws.Range("H73").Value = 5427.375
$ws.Range("J73").Value = 3651.842
$ws.Range("L73").Value = 10955.526
$ws.Range("N73").Value = -12827.526

$ws.Range("H135").Value = 5000500
$ws.Range("I135").Value = 10000000
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 90000000
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -89997465
$ws.Range("N135").Value = -14070

$ws.Range("H138").Value = 2226772.5
$ws.Range("I138").Value = 1758.5454
$ws.Range("J138").Value = 4355046.5
$ws.Range("K138").Value = 5275.6362
$ws.Range("L138").Value = 13065139.5
$ws.Range("M138").Value = -135.6361999999999
$ws.Range("N138").Value = -13075419.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2824.7778
$ws.Range("I45").Value = 2674.4
$ws.Range("K45").Value = 2674.4
$ws.Range("M45").Value = -2297.4

$ws.Range("H63").Value = 1700
$ws.Range("I63").Value = 1400
$ws.Range("K63").Value = 1400
$ws.Range("M63").Value = -714

$ws.Range("H66").Value = 1700
$ws.Range("I66").Value = 1400
$ws.Range("K66").Value = 7000
$ws.Range("M66").Value = -3568

$ws.Range("H74").Value = 48839
$ws.Range("I74").Value = 78764.92
$ws.Range("K74").Value = 78764.92
$ws.Range("M74").Value = -77890.92

$ws.Range("H77").Value = 48839
$ws.Range("I77").Value = 78764.92
$ws.Range("K77").Value = 393824.6
$ws.Range("M77").Value = -389456.6

$ws.Range("H97").Value = 2526125
$ws.Range("I97").Value = 612.9655
$ws.Range("K97").Value = 612.9655
$ws.Range("M97").Value = -116.9655

$ws.Range("H132").Value = 5452.375
$ws.Range("I132").Value = 2907.5483
$ws.Range("K132").Value = 8722.644899999999
$ws.Range("M132").Value = -6192.644899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5378330.5
$ws.Range("I20").Value = 7247987.5
$ws.Range("J20").Value = 3066.5
$ws.Range("K20").Value = 7247987.5
$ws.Range("L20").Value = 3066.5
$ws.Range("M20").Value = -7247740.5
$ws.Range("N20").Value = -3560.5

$ws.Range("H22").Value = 7936757.5
$ws.Range("I22").Value = 7936757.5
$ws.Range("K22").Value = 7936757.5
$ws.Range("M22").Value = -7936584.5

$ws.Range("H80").Value = 38462084
$ws.Range("J80").Value = 599.8
$ws.Range("L80").Value = 599.8
$ws.Range("N80").Value = -2595.8

$ws.Range("H83").Value = 38462084
$ws.Range("J83").Value = 599.8
$ws.Range("L83").Value = 2999
$ws.Range("N83").Value = -12983

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7531.1465
$ws.Range("I31").Value = 2631.5789
$ws.Range("J31").Value = 11762.591
$ws.Range("K31").Value = 2631.5789
$ws.Range("L31").Value = 11762.591
$ws.Range("M31").Value = -2336.5789
$ws.Range("N31").Value = -12352.591

$ws.Range("H34").Value = 7531.1465
$ws.Range("I34").Value = 2631.5789
$ws.Range("J34").Value = 11762.591
$ws.Range("K34").Value = 2631.5789
$ws.Range("L34").Value = 11762.591
$ws.Range("M34").Value = -2429.5789
$ws.Range("N34").Value = -12166.591

$ws.Range("H86").Value = 64956980
$ws.Range("I86").Value = 30337632
$ws.Range("J86").Value = 111116110
$ws.Range("K86").Value = 30337632
$ws.Range("L86").Value = 111116110
$ws.Range("M86").Value = -30336509
$ws.Range("N86").Value = -111118356

$ws.Range("H89").Value = 64956980
$ws.Range("I89").Value = 30337632
$ws.Range("J89").Value = 111116110
$ws.Range("K89").Value = 151688160
$ws.Range("L89").Value = 555580550
$ws.Range("M89").Value = -151682544
$ws.Range("N89").Value = -555591782

$ws.Range("H94").Value = 1405.7858
$ws.Range("I94").Value = 2433
$ws.Range("K94").Value = 2433
$ws.Range("M94").Value = -1982

$ws.Range("H99").Value = 4195.5835
$ws.Range("I99").Value = 3125.8667
$ws.Range("J99").Value = 5978.4443
$ws.Range("K99").Value = 3125.8667
$ws.Range("L99").Value = 5978.4443
$ws.Range("M99").Value = -1627.8667
$ws.Range("N99").Value = -8974.444299999999

$ws.Range("H105").Value = 10204695
$ws.Range("I105").Value = 14286234
$ws.Range("J105").Value = 850
$ws.Range("K105").Value = 14286234
$ws.Range("L105").Value = 850
$ws.Range("M105").Value = -14284487
$ws.Range("N105").Value = -4344

$ws.Range("H122").Value = 3833.8064
$ws.Range("I122").Value = 2405.111
$ws.Range("K122").Value = 7215.333
$ws.Range("M122").Value = -4765.333

$ws.Range("H126").Value = 4195.5835
$ws.Range("I126").Value = 3125.8667
$ws.Range("J126").Value = 5978.4443
$ws.Range("K126").Value = 9377.6001
$ws.Range("L126").Value = 17935.3329
$ws.Range("M126").Value = -6907.6001
$ws.Range("N126").Value = -22875.3329

$ws.Range("H132").Value = 4849.154
$ws.Range("I132").Value = 2536.6191
$ws.Range("K132").Value = 7609.8573
$ws.Range("M132").Value = -5079.8573

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J38").Value = 125002510
$ws.Range("L38").Value = 375007530
$ws.Range("N38").Value = -375008224

$ws.Range("H40").Value = 215.33333
$ws.Range("J40").Value = 297.75
$ws.Range("L40").Value = 1191
$ws.Range("N40").Value = -1329

$ws.Range("H97").Value = 538.7778
$ws.Range("J97").Value = 640.6
$ws.Range("L97").Value = 1921.8
$ws.Range("N97").Value = -2913.8

$ws.Range("H113").Value = 8448.200000000001
$ws.Range("J113").Value = 9275.777
$ws.Range("L113").Value = 27827.331
$ws.Range("N113").Value = -32167.331

$ws.Range("H132").Value = 6090
$ws.Range("J132").Value = 9846.532999999999
$ws.Range("L132").Value = 88618.79699999999
$ws.Range("N132").Value = -93678.79699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7283093.5
$ws.Range("I46").Value = 2654981.5
$ws.Range("K46").Value = 2654981.5
$ws.Range("M46").Value = -2654793.5

$ws.Range("H68").Value = 166672660
$ws.Range("I68").Value = 333337000
$ws.Range("K68").Value = 333337000
$ws.Range("M68").Value = -333336251

$ws.Range("H71").Value = 166672660
$ws.Range("I71").Value = 333337000
$ws.Range("K71").Value = 1666685000
$ws.Range("M71").Value = -1666681256

$ws.Range("H107").Value = 3539.6
$ws.Range("I107").Value = 3539.6
$ws.Range("K107").Value = 3539.6
$ws.Range("M107").Value = -1619.6

$ws.Range("H132").Value = 6136.88
$ws.Range("I132").Value = 3199.7727
$ws.Range("K132").Value = 9599.3181
$ws.Range("M132").Value = -7069.3181

$ws.Range("H136").Value = 14352.676
$ws.Range("I136").Value = 3539.2856
$ws.Range("J136").Value = 20934.738
$ws.Range("K136").Value = 10617.8568
$ws.Range("L136").Value = 62804.21400000001
$ws.Range("M136").Value = -8067.856800000001
$ws.Range("N136").Value = -67904.21400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 68000
$ws.Range("I47").Value = 68000
$ws.Range("K47").Value = 68000
$ws.Range("M47").Value = -67428

$ws.Range("H49").Value = 208000
$ws.Range("I49").Value = 208000
$ws.Range("K49").Value = 208000
$ws.Range("M49").Value = -207770

$ws.Range("H62").Value = 8705.5
$ws.Range("I62").Value = 8705.5
$ws.Range("K62").Value = 8705.5
$ws.Range("M62").Value = -8081.5

$ws.Range("H65").Value = 8705.5
$ws.Range("I65").Value = 8705.5
$ws.Range("K65").Value = 43527.5
$ws.Range("M65").Value = -40407.5

$ws.Range("H132").Value = 55563944
$ws.Range("I132").Value = 83342830
$ws.Range("K132").Value = 250028490
$ws.Range("M132").Value = -250025960

$ws.Range("H136").Value = 421154.62
$ws.Range("I136").Value = 1284.6923
$ws.Range("J136").Value = 917364.5600000001
$ws.Range("K136").Value = 3854.0769
$ws.Range("L136").Value = 2752093.68
$ws.Range("M136").Value = -1304.0769
$ws.Range("N136").Value = -2757193.68
